# Version 2 de ExtractPDF
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# --------------------------------------------------------------------
# 1) Update existing "cotizacion" value (D2)
# --------------------------------------------------------------------
$ws.Range("D2").Value = 20654

# --------------------------------------------------------------------
# 2) Write all the new cell VALUES first (H1:K17)
# --------------------------------------------------------------------
$ws.Range("H1").Value = "Sequence #"
$ws.Range("I1").Value = "Upper Bound"
$ws.Range("J1").Value = "Rate"
$ws.Range("K1").Value = "Factor"

$seq    = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
$upper  = @(500,1000,2000,3000,4000,5000,6000,7000,8000,9000,10000,11000,12000,13000,20000,9999999)
$rate   = @(0.35,0.33,0.29,0.26,0.22,0.21,0.2,0.19,0.18,0.17,0.155,0.14,0.1375,0.135,0.13,0.125)
$factor = @(0.1279352436098831,0.126940954842813,0.1249641727581296,0.1234919553361154,0.1215429015542623,0.1210581306386495,0.1205743596691218,0.1200915900996986,0.1196098233768364,0.1191290609393636,0.1184098034467487,0.1176928136117802,0.117573536107308,0.1174543217474871,0.117216082549313,0.116978096191973)

for ($i = 0; $i -lt $seq.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 8).Value  = $seq[$i]
    $ws.Cells.Item($r, 9).Value  = $upper[$i]
    $ws.Cells.Item($r, 10).Value = $rate[$i]
    $ws.Cells.Item($r, 11).Value = $factor[$i]
}

# --------------------------------------------------------------------
# 3) Establish each distinct style on ONE representative cell, in the
#    exact order the styles first appear, then fan it out with
#    Copy + PasteSpecial(Formats) so no stray intermediate xf entries
#    are left behind in the stylesheet.
# --------------------------------------------------------------------

# -- style: fontId=1 (Calibri, no family) plain -> general "H/I/K" cells
$ws.Cells.Item(2, 8).Font.Name = "Calibri"
$ws.Cells.Item(2, 8).Copy()
$ws.Range("H3:H17").PasteSpecial($xlPasteFormats)
$ws.Range("I2").PasteSpecial($xlPasteFormats)
$ws.Range("K2").PasteSpecial($xlPasteFormats)
$ws.Range("K4:K17").PasteSpecial($xlPasteFormats)

# -- style: fontId=2 (Aptos Narrow) + numFmtId=10 (0.00%) -> Rate column
$ws.Cells.Item(2, 10).Font.Name = "Aptos Narrow"
$ws.Cells.Item(2, 10).NumberFormat = "0.00%"
$ws.Cells.Item(2, 10).Copy()
$ws.Range("J3:J17").PasteSpecial($xlPasteFormats)

# -- style: fontId=3 (Arial) plain -> K3 only
$ws.Cells.Item(3, 11).Font.Name = "Arial"

# -- style: fontId=1 (Calibri, no family) + center/center alignment -> headers
$h1 = $ws.Cells.Item(1, 8)
$h1.Font.Name = "Calibri"
$h1.VerticalAlignment = $xlCenter
$h1.HorizontalAlignment = $xlCenter
$h1.Copy()
$ws.Range("I1:K1").PasteSpecial($xlPasteFormats)

# -- style: fontId=2 (Aptos Narrow) + numFmtId=3 (#,##0) -> Upper Bound column (rows 3-17)
$ws.Cells.Item(3, 9).Font.Name = "Aptos Narrow"
$ws.Cells.Item(3, 9).NumberFormat = "#,##0"
$ws.Cells.Item(3, 9).Copy()
$ws.Range("I4:I17").PasteSpecial($xlPasteFormats)

# --------------------------------------------------------------------
# 4) Column widths (best-fit H:I)
# --------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 10.25
$ws.Columns.Item(9).ColumnWidth = 12.2

# --------------------------------------------------------------------
# 5) Selection
# --------------------------------------------------------------------
$ws.Range("K2").Select()
